$wb = $excel.ActiveWorkbook
$wsPurchase = $wb.Worksheets.Item("Purchase 22-23")

# 1. Break the F29 formula on "Purchase 22-23" (typo "G29E29") so it
#    now evaluates to a #NAME? error instead of -57151.
$wsPurchase.Range("F29").Formula = "=F28+G29E29"

# 2. Delete the blank row 36 on "Purchase 22-23". This shifts the data
#    that used to live on rows 37 and 39 up by one row (-> rows 36, 38),
#    shrinking the used range from A1:F39 down to A1:F38.
$wsPurchase.Rows("36").Delete()

# 3. "Purchase 22-23" becomes the active/selected sheet (previously it
#    was "Sale 22-23"), with F29 as the selected cell.
$wsPurchase.Activate()
$wsPurchase.Range("F29").Select()
